$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 174.14285
$ws.Range("I12").Value = 163.33333
$ws.Range("J12").Value = 182.25
$ws.Range("K12").Value = 163.33333
$ws.Range("L12").Value = 182.25
$ws.Range("M12").Value = 6.666670000000011
$ws.Range("N12").Value = -522.25

$ws.Range("H28").Value = 212.53847
$ws.Range("I28").Value = 229.1
$ws.Range("J28").Value = 157.33333
$ws.Range("K28").Value = 229.1
$ws.Range("L28").Value = 157.33333
$ws.Range("M28").Value = 255.9
$ws.Range("N28").Value = -1127.33333

$ws.Range("H31").Value = 3133
$ws.Range("I31").Value = 1718.1111
$ws.Range("J31").Value = 9500
$ws.Range("K31").Value = 5154.3333
$ws.Range("L31").Value = 28500
$ws.Range("M31").Value = -4924.3333
$ws.Range("N31").Value = -28960

$ws.Range("H58").Value = 666.46155
$ws.Range("I58").Value = 366.4
$ws.Range("J58").Value = 1666.6666
$ws.Range("K58").Value = 1099.2
$ws.Range("L58").Value = 4999.9998
$ws.Range("M58").Value = -949.1999999999998
$ws.Range("N58").Value = -5299.9998

$ws.Range("H131").Value = 1058.65
$ws.Range("I131").Value = 723.3125
$ws.Range("J131").Value = 2400
$ws.Range("K131").Value = 2169.9375
$ws.Range("L131").Value = 7200
$ws.Range("M131").Value = 2870.0625
$ws.Range("N131").Value = -17280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 16564.428
$ws.Range("I36").Value = 2150.4
$ws.Range("J36").Value = 52599.5
$ws.Range("K36").Value = 2150.4
$ws.Range("L36").Value = 52599.5
$ws.Range("M36").Value = -1804.4
$ws.Range("N36").Value = -53291.5

$ws.Range("H61").Value = 3601.0303
$ws.Range("I61").Value = 2014.75
$ws.Range("J61").Value = 5094
$ws.Range("K61").Value = 2014.75
$ws.Range("L61").Value = 5094
$ws.Range("M61").Value = -1802.75
$ws.Range("N61").Value = -5518

$ws.Range("H97").Value = 1283.1818
$ws.Range("I97").Value = 1283.1818
$ws.Range("K97").Value = 1283.1818
$ws.Range("M97").Value = -787.1818000000001

$ws.Range("H110").Value = 930.2857
$ws.Range("I110").Value = 793.0909
$ws.Range("J110").Value = 1433.3334
$ws.Range("K110").Value = 793.0909
$ws.Range("L110").Value = 1433.3334
$ws.Range("M110").Value = 1251.9091
$ws.Range("N110").Value = -5523.3334

$ws.Range("H136").Value = 3601.0303
$ws.Range("I136").Value = 2014.75
$ws.Range("J136").Value = 5094
$ws.Range("K136").Value = 6044.25
$ws.Range("L136").Value = 15282
$ws.Range("M136").Value = -3494.25
$ws.Range("N136").Value = -20382

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 304
$ws.Range("I11").Value = 257.66666
$ws.Range("J11").Value = 338.75
$ws.Range("K11").Value = 257.66666
$ws.Range("L11").Value = 338.75
$ws.Range("M11").Value = -117.66666
$ws.Range("N11").Value = -618.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 825.8125
$ws.Range("I16").Value = 800.9286
$ws.Range("K16").Value = 800.9286
$ws.Range("M16").Value = -513.9286

$ws.Range("H105").Value = 1452.8182
$ws.Range("I105").Value = 1692.625
$ws.Range("J105").Value = 813.3333
$ws.Range("K105").Value = 1692.625
$ws.Range("L105").Value = 813.3333
$ws.Range("M105").Value = 54.375
$ws.Range("N105").Value = -4307.3333

$ws.Range("H113").Value = 825.8125
$ws.Range("I113").Value = 800.9286
$ws.Range("K113").Value = 800.9286
$ws.Range("M113").Value = 1369.0714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 357.14285
$ws.Range("I8").Value = 357.14285
$ws.Range("K8").Value = 1071.42855
$ws.Range("M8").Value = -932.4285500000001

$ws.Range("H23").Value = 61.333332
$ws.Range("I23").Value = 16.4
$ws.Range("J23").Value = 117.5
$ws.Range("K23").Value = 49.2
$ws.Range("L23").Value = 352.5
$ws.Range("M23").Value = 185.8
$ws.Range("N23").Value = -822.5

$ws.Range("H44").Value = 935.8
$ws.Range("I44").Value = 191
$ws.Range("J44").Value = 1432.3334
$ws.Range("K44").Value = 573
$ws.Range("L44").Value = 4297.0002
$ws.Range("M44").Value = -175
$ws.Range("N44").Value = -5093.0002

$ws.Range("H113").Value = 1642
$ws.Range("I113").Value = 2529.4285
$ws.Range("J113").Value = 754.5714
$ws.Range("K113").Value = 7588.2855
$ws.Range("L113").Value = 2263.7142
$ws.Range("M113").Value = -5418.2855
$ws.Range("N113").Value = -6603.7142

$ws.Range("H131").Value = 935.4
$ws.Range("I131").Value = 605
$ws.Range("J131").Value = 978.0323
$ws.Range("K131").Value = 1815
$ws.Range("L131").Value = 2934.0969
$ws.Range("M131").Value = 3225
$ws.Range("N131").Value = -13014.0969

$ws.Range("H132").Value = 1470.2
$ws.Range("I132").Value = 931.2857
$ws.Range("J132").Value = 1760.3846
$ws.Range("K132").Value = 8381.5713
$ws.Range("L132").Value = 15843.4614
$ws.Range("M132").Value = -5851.5713
$ws.Range("N132").Value = -20903.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 1833.3334
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1833.3334
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -211.3334
$ws.Range("N102").Value = -5744

$ws.Range("H122").Value = 4943.5
$ws.Range("I122").Value = 1676
$ws.Range("J122").Value = 6428.727
$ws.Range("K122").Value = 5028
$ws.Range("L122").Value = 19286.181
$ws.Range("M122").Value = -2578
$ws.Range("N122").Value = -24186.181

$ws.Range("H132").Value = 2020.7222
$ws.Range("I132").Value = 1552.4348
$ws.Range("K132").Value = 4657.3044
$ws.Range("M132").Value = -2127.3044

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 23500
$ws.Range("J41").Value = 12000
$ws.Range("L41").Value = 12000
$ws.Range("N41").Value = -12780

$ws.Range("H47").Value = 70000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 70000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 70000
$ws.Range("N47").Value = -71144
$ws.Range("M47").ClearContents()

$ws.Range("H54").Value = 23146.2
$ws.Range("I54").Value = 13000
$ws.Range("J54").Value = 24273.555
$ws.Range("K54").Value = 13000
$ws.Range("L54").Value = 24273.555
$ws.Range("M54").Value = -12480
$ws.Range("N54").Value = -25313.555

$ws.Range("H81").Value = 888.1212
$ws.Range("I81").Value = 747.2
$ws.Range("J81").Value = 1328.5
$ws.Range("K81").Value = 1494.4
$ws.Range("L81").Value = 2657
$ws.Range("M81").Value = -433.4000000000001
$ws.Range("N81").Value = -4779

$ws.Range("H84").Value = 888.1212
$ws.Range("I84").Value = 747.2
$ws.Range("J84").Value = 1328.5
$ws.Range("K84").Value = 7472
$ws.Range("L84").Value = 13285
$ws.Range("M84").Value = -2168
$ws.Range("N84").Value = -23893
